# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 34
    "F3"  = 6280
    "F4"  = 176
    "F7"  = 1893
    "F8"  = 1427
    "F9"  = 296
    "F10" = 956
    "F11" = 250
    "F12" = 5584
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
